$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "id"
$ws.Range("B1").Value = "question"
$ws.Range("C1").Value = "option1"
$ws.Range("D1").Value = "option2"
$ws.Range("E1").Value = "option3"
$ws.Range("F1").Value = "option4"
$ws.Range("G1").Value = "answer"
